$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-18.851413554346, 1.428320567775983, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346),
    @(-18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, 1.113216394376549, -18.851413554346, -18.851413554346),
    @(-18.851413554346, 1.577323838568174, 4.321925197159212, -18.851413554346, 3.390012508284986, -18.851413554346, 1.471307919094931, -18.851413554346, 0.7904353440846377, -18.851413554346),
    @(-18.851413554346, 2.042179056810245, -18.851413554346, -18.851413554346, -18.851413554346, 2.81926041373017, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346),
    @(-18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346),
    @(2.533889080975926, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346),
    @(-18.851413554346, -18.851413554346, -18.851413554346, 1.198542650929214, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346),
    @(3.828698041319877, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346),
    @(-18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, 1.75043429439709, -18.851413554346, 2.23539973421075),
    @(-18.851413554346, -18.851413554346, -18.851413554346, 3.009969110778929, -18.851413554346, 2.84103210536718, -18.851413554346, -18.851413554346, -18.851413554346, 1.945314483202426),
    @(-18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346),
    @(-18.851413554346, -18.851413554346, -18.851413554346, 2.691422312260855, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, 1.666637315478572, 1.759795899164231),
    @(-18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, 1.913086335188555),
    @(-18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346),
    @(-18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, 1.910427378750442, -18.851413554346),
    @(-18.851413554346, 2.451580613806552, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, 2.056402057175585, 2.09743629540746, 2.560358601415702, -18.851413554346),
    @(-18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, 2.01915661254856, 2.0634540209277, 2.443046218567352, -18.851413554346),
    @(-18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, -18.851413554346, 1.641086112992164, 1.863151442226432, -18.851413554346, -18.851413554346),
    @(-18.851413554346, 1.277442961208757, -18.851413554346, -18.851413554346, 3.250465004324747, -18.851413554346, 1.663251385507641, 1.247771140183432, -18.851413554346, 2.1000284060384),
    @(-18.851413554346, 1.208872054704775, -18.851413554346, 1.673519304326396, -18.851413554346, 2.530229197987196, 1.447840307881874, -18.851413554346, -18.851413554346, -18.851413554346)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $c = $j + 2
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}
